$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Application Tracker")

# Existing row 4 (Pretend Corp. / Software Developer): status flips from Open to Closed
$ws.Range("L4").Value = "Closed"

# Row 5: Bullhorn / Junior Software Engineer / not yet / Aqeelah Jones / "" / Ready / URL / Open / JavaScript
$ws.Range("B5").Value = "Bullhorn"
$ws.Range("C5").Value = "Junior Software Engineer"
$ws.Range("D5").Value = "not yet"
$ws.Range("F5").Value = "Aqeelah Jones"
$ws.Range("G5").Value = """"""
$ws.Range("J5").Value = "Ready"
$ws.Range("K5").Value = "https://www.linkedin.com/jobs/search/?currentJobId=890579596&keywords=software%20developer&location=Greater%20Atlanta%20Area&locationId=us%3A52"
$ws.Range("L5").Value = "Open"
$ws.Range("M5").Value = "JavaScript"

# Row 6: Infor / Junior Software Engineer / not yet / Allie Persinger / "" / Ready / URL / Open / JavaScript
$ws.Range("B6").Value = "Infor"
$ws.Range("C6").Value = "Junior Software Engineer"
$ws.Range("D6").Value = "not yet"
$ws.Range("F6").Value = "Allie Persinger"
$ws.Range("G6").Value = """"""
$ws.Range("J6").Value = "Ready"
$ws.Range("K6").Value = "https://www.linkedin.com/jobs/search/?currentJobId=889095712&keywords=software%20developer&location=Greater%20Atlanta%20Area&locationId=us%3A52"
$ws.Range("L6").Value = "Open"
$ws.Range("M6").Value = "JavaScript"

# Row 7: Brooksource / Junior Software Engineer / not yet / Wynne Rosenbleeth / "" / 470-419-2504 / Ready / URL / Open / JavaScript
$ws.Range("B7").Value = "Brooksource"
$ws.Range("C7").Value = "Junior Software Engineer"
$ws.Range("D7").Value = "not yet"
$ws.Range("F7").Value = "Wynne Rosenbleeth"
$ws.Range("G7").Value = """"""
$ws.Range("H7").Value = "470-419-2504"
$ws.Range("J7").Value = "Ready"
$ws.Range("K7").Value = "https://www.linkedin.com/jobs/search/?currentJobId=888000630&keywords=software%20developer&location=Greater%20Atlanta%20Area&locationId=us%3A52"
$ws.Range("L7").Value = "Open"
$ws.Range("M7").Value = "JavaScript"

$wb.Save()
